$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-NumericTextCell($addr, $val) {
    # Excel auto-detects digit-only strings as numbers, which would lose
    # the trailing newline and store the cell as a numeric value instead
    # of a shared string. Briefly force Text format, isolated to this one
    # cell only (its own WrapText toggle keeps it from sharing a style
    # slot with any other cell while the number format is mutated), so
    # the value lands as literal text; then restore General format and
    # wrap state so the only difference left is the real final style.
    $ws.Range($addr).WrapText = $true
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = $null
    $ws.Range($addr).WrapText = $false
}

Set-TextCell "A2" "Dynazzy`n"
Set-NumericTextCell "B2" "279`n"
Set-NumericTextCell "C2" "2170`n"
Set-TextCell "D2" "1 725,00 `n"

Set-TextCell "A3" "Yakidoo`n"
Set-NumericTextCell "B3" "609`n"
Set-NumericTextCell "C3" "3031`n"
Set-TextCell "D3" "650,00 `n"

Set-TextCell "A4" "Jamia`n"
Set-NumericTextCell "B4" "97`n"
Set-NumericTextCell "C4" "5536`n"
Set-TextCell "D4" "925,00 `n"

Set-TextCell "A5" "Mydo`n"
Set-NumericTextCell "B5" "570`n"
Set-NumericTextCell "C5" "2223`n"
Set-TextCell "D5" "675,00 `n"

$ws.Range("A2:D5").WrapText = $true
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
